$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: copy header style from A1 to S1:AL1, then set text ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("S1:AL1").PasteSpecial(-4122) | Out-Null

$ws.Range("S1").Value = "question_3_hint_used"
$ws.Range("T1").Value = "question_3_tries_taken"
$ws.Range("U1").Value = "question_3_correct"
$ws.Range("V1").Value = "question_3_selections"
$ws.Range("W1").Value = "question_4_hint_used"
$ws.Range("X1").Value = "question_4_tries_taken"
$ws.Range("Y1").Value = "question_4_correct"
$ws.Range("Z1").Value = "question_4_selections"
$ws.Range("AA1").Value = "question_5_hint_used"
$ws.Range("AB1").Value = "question_5_tries_taken"
$ws.Range("AC1").Value = "question_5_correct"
$ws.Range("AD1").Value = "question_5_selections"
$ws.Range("AE1").Value = "question_6_hint_used"
$ws.Range("AF1").Value = "question_6_tries_taken"
$ws.Range("AG1").Value = "question_6_correct"
$ws.Range("AH1").Value = "question_6_selections"
$ws.Range("AI1").Value = "question_7_hint_used"
$ws.Range("AJ1").Value = "question_7_tries_taken"
$ws.Range("AK1").Value = "question_7_correct"
$ws.Range("AL1").Value = "question_7_selections"

# --- Row 3: new data row for "Autonomous Vehicle" ---
# Copy numeric date style from B2:C2 to B3:C3 first
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B3:C3").PasteSpecial(-4122) | Out-Null

$ws.Range("A3").Value = "Autonomous Vehicle"
$ws.Range("B3").Value = 45873.12718336806
$ws.Range("C3").Value = 45873.12720300926
$ws.Range("D3").Value = 1697
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = $false
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = $false
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = $false
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = $false
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = $false
$ws.Range("R3").Value = ""
$ws.Range("S3").Value = $false
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = $false
$ws.Range("V3").Value = ""
$ws.Range("W3").Value = $false
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = $false
$ws.Range("Z3").Value = ""
$ws.Range("AA3").Value = $false
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = $false
$ws.Range("AD3").Value = ""
$ws.Range("AE3").Value = $false
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = $false
$ws.Range("AH3").Value = ""
$ws.Range("AI3").Value = $false
$ws.Range("AJ3").Value = 0
$ws.Range("AK3").Value = $false
$ws.Range("AL3").Value = ""
